$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
